$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of random-method stats appended below the existing data (row 4 -> row 5).
# Copy the date/time formatting from the cell above so the new date cell reuses
# the existing style instead of minting a new one.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A5").Value = 42602.016458333332
$ws.Range("B5").Value = "Random"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 11
$ws.Range("I5").Value = 89
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 60
$ws.Range("M5").Value = 40
